# Update crypto price/volume cells per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.438.82"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.697.24"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.88"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5480"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2735"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06445"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.98"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07687"
$ws.Range("E11").Value = "  +2.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.702.96"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.552"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008407"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.71"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.490.24"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.944"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.99"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.25"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.256"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.94"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("E25").Value = "  +6.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.896"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.80"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06259"
$ws.Range("E28").Value = "  -5.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.380"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.333"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.614"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.596"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.690"
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.039"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6172"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.764"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.117.39"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.130"
$ws.Range("E40").Value = "  -3.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8796"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.20"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.849.84"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.58"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05287"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.127"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4302"
$ws.Range("E51").Value = "  +0.02%  "
